# Regenerate s_vals data to filter save games.
# Updates the numeric columns B:G for rows 2-29 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054),
    @(0.06328177979961902, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1, 1.328214039578707),
    @(0.1554434735375247, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.145820798638228),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(0.3464964993005633, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 3.21235251628842),
    @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 0, 4.266567469086664),
    @(0.001754667048134761, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 0, 0.7068598578304137),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(0.3464964993005633, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 2.652525447291612),
    @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 5.964442013611383),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126),
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729),
    @(0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 6.48142807727062, 1, 7.878242263594639),
    @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 6.48142807727062, 0, 11.94598338380795),
    @(0.001754667048134761, 0.0001537489499301437, 0.1529057820181812, 0.4998867070740569, 0, 0.654700905090303),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126),
    @(1.505614041169197, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 20.64246832346449),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(0.7287194209349384, 1.65323645889881, 157.8057217802531, 6.48142807727062, 1, 166.6691057373575),
    @(0.02258322285507441, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.012960547955778),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126),
    @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126),
    @(3.182878228561681, 9.226618575922256, 3.082599426703578, 6.48142807727062, 1, 21.97352430845813),
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = 2 + $j   # column B = 2
        $ws.Cells.Item($r, $col).Value = $rowValues[$j]
    }
}
